$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 135, shifting existing rows 135:180 down to 136:181
$ws.Rows("135:135").Insert()

$ws.Range("A135").Value = 4
$ws.Range("B135").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C135").Value = "Los Lagos"
$ws.Range("D135").Value = 44524
$ws.Range("E135").Value = 10
$ws.Range("F135").Value = "Fruta"
$ws.Range("G135").Value = 100102
$ws.Range("H135").Value = "Cítricos"
$ws.Range("I135").Value = 100102006
$ws.Range("J135").Value = "Pomelo"
$ws.Range("K135").Value = "Start Ruby"
$ws.Range("L135").Value = "Primera"
$ws.Range("M135").Value = 160
$ws.Range("N135").Value = 11000
$ws.Range("O135").Value = 12000
$ws.Range("P135").Value = 11500
$ws.Range("Q135").Value = "`$/caja 14 kilos empedrada"
$ws.Range("R135").Value = "Región de O'Higgins"
$ws.Range("S135").Value = 821
$ws.Range("T135").Value = 14
